$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H12").Value = 96
$ws_ALC.Range("I12").Value = 75
$ws_ALC.Range("J12").Value = 127.5
$ws_ALC.Range("K12").Value = 75
$ws_ALC.Range("L12").Value = 127.5
$ws_ALC.Range("M12").Value = 95
$ws_ALC.Range("N12").Value = -467.5
$ws_ALC.Range("H19").Value = 468.33334
$ws_ALC.Range("J19").Value = 468.33334
$ws_ALC.Range("L19").Value = 468.33334
$ws_ALC.Range("N19").Value = -818.33334
$ws_ALC.Range("H40").Value = 5809.9
$ws_ALC.Range("I40").Value = 4014.1428
$ws_ALC.Range("K40").Value = 4014.1428
$ws_ALC.Range("M40").Value = -3839.1428
$ws_ALC.Range("H41").Value = 617
$ws_ALC.Range("J41").Value = 624
$ws_ALC.Range("L41").Value = 624
$ws_ALC.Range("N41").Value = -1504
$ws_ALC.Range("H53").Value = 325.75
$ws_ALC.Range("I53").Value = 358.66666
$ws_ALC.Range("K53").Value = 358.66666
$ws_ALC.Range("M53").Value = 278.33334
$ws_ALC.Range("H88").Value = 2763.8572
$ws_ALC.Range("J88").Value = 2986.6
$ws_ALC.Range("L88").Value = 2986.6
$ws_ALC.Range("N88").Value = -3798.6
$ws_ALC.Range("H91").Value = 2763.8572
$ws_ALC.Range("J91").Value = 2986.6
$ws_ALC.Range("L91").Value = 2986.6
$ws_ALC.Range("N91").Value = -5794.6
$ws_ALC.Range("H96").Value = 1615.75
$ws_ALC.Range("I96").Value = 418.14285
$ws_ALC.Range("K96").Value = 1254.42855
$ws_ALC.Range("M96").Value = 118.5714499999999
$ws_ALC.Range("H98").Value = 2357.5386
$ws_ALC.Range("I98").Value = 850.9091
$ws_ALC.Range("K98").Value = 850.9091
$ws_ALC.Range("M98").Value = 647.0909
$ws_ALC.Range("H107").Value = 728.3333
$ws_ALC.Range("I107").Value = 653.5294
$ws_ALC.Range("K107").Value = 653.5294
$ws_ALC.Range("M107").Value = 1266.4706
$ws_ALC.Range("H122").Value = 2357.5386
$ws_ALC.Range("I122").Value = 850.9091
$ws_ALC.Range("K122").Value = 2552.7273
$ws_ALC.Range("M122").Value = -102.7273
$ws_ALC.Range("H127").Value = 777
$ws_ALC.Range("I127").Value = 0
$ws_ALC.Range("K127").Value = 0
$ws_ALC.Range("M127").ClearContents()
$ws_ALC.Range("H137").Value = 2641
$ws_ALC.Range("I137").Value = 2469.4
$ws_ALC.Range("J137").Value = 3499
$ws_ALC.Range("K137").Value = 7408.200000000001
$ws_ALC.Range("L137").Value = 10497
$ws_ALC.Range("M137").Value = -4858.200000000001
$ws_ALC.Range("N137").Value = -15597
$ws_ALC.Range("H138").Value = 6489.353
$ws_ALC.Range("I138").Value = 4500
$ws_ALC.Range("J138").Value = 6754.6
$ws_ALC.Range("K138").Value = 13500
$ws_ALC.Range("L138").Value = 20263.8
$ws_ALC.Range("M138").Value = -8360
$ws_ALC.Range("N138").Value = -30543.8
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H74").Value = 21004.166
$ws_ARM.Range("I74").Value = 16501
$ws_ARM.Range("K74").Value = 16501
$ws_ARM.Range("M74").Value = -15627
$ws_ARM.Range("H77").Value = 21004.166
$ws_ARM.Range("I77").Value = 16501
$ws_ARM.Range("K77").Value = 82505
$ws_ARM.Range("M77").Value = -78137
$ws_ARM.Range("H110").Value = 3763.5
$ws_ARM.Range("I110").Value = 1405.5555
$ws_ARM.Range("J110").Value = 10837.333
$ws_ARM.Range("K110").Value = 1405.5555
$ws_ARM.Range("L110").Value = 10837.333
$ws_ARM.Range("M110").Value = 639.4445000000001
$ws_ARM.Range("N110").Value = -14927.333
$ws_ARM.Range("H122").Value = 2000
$ws_ARM.Range("I122").Value = 2000
$ws_ARM.Range("K122").Value = 6000
$ws_ARM.Range("M122").Value = -3550
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H9").Value = 65000
$ws_BSM.Range("J9").Value = 65000
$ws_BSM.Range("L9").Value = 65000
$ws_BSM.Range("N9").Value = -65336
$ws_BSM.Range("H22").Value = 0
$ws_BSM.Range("I22").Value = 0
$ws_BSM.Range("K22").Value = 0
$ws_BSM.Range("M22").ClearContents()
$ws_BSM.Range("H86").Value = 11455.333
$ws_BSM.Range("J86").Value = 25000
$ws_BSM.Range("L86").Value = 25000
$ws_BSM.Range("N86").Value = -27246
$ws_BSM.Range("H89").Value = 11455.333
$ws_BSM.Range("J89").Value = 25000
$ws_BSM.Range("L89").Value = 125000
$ws_BSM.Range("N89").Value = -136232
$ws_BSM.Range("H134").Value = 1102.2
$ws_BSM.Range("I134").Value = 1102.2
$ws_BSM.Range("K134").Value = 3306.6
$ws_BSM.Range("M134").Value = -771.6000000000004
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H99").Value = 3909.5
$ws_CRP.Range("I99").Value = 3912
$ws_CRP.Range("K99").Value = 3912
$ws_CRP.Range("M99").Value = -2414
$ws_CRP.Range("H126").Value = 3909.5
$ws_CRP.Range("I126").Value = 3912
$ws_CRP.Range("K126").Value = 11736
$ws_CRP.Range("M126").Value = -9266
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H75").Value = 3998
$ws_CUL.Range("J75").Value = 3998
$ws_CUL.Range("L75").Value = 11994
$ws_CUL.Range("N75").Value = -13990
$ws_CUL.Range("H78").Value = 3998
$ws_CUL.Range("J78").Value = 3998
$ws_CUL.Range("L78").Value = 35982
$ws_CUL.Range("N78").Value = -45966
$ws_CUL.Range("H97").Value = 1290.1666
$ws_CUL.Range("J97").Value = 793.5
$ws_CUL.Range("L97").Value = 2380.5
$ws_CUL.Range("N97").Value = -3372.5
$ws_CUL.Range("H98").Value = 2962.7144
$ws_CUL.Range("I98").Value = 3444.75
$ws_CUL.Range("J98").Value = 2320
$ws_CUL.Range("K98").Value = 10334.25
$ws_CUL.Range("L98").Value = 6960
$ws_CUL.Range("M98").Value = -8836.25
$ws_CUL.Range("N98").Value = -9956
$ws_CUL.Range("H109").Value = 4301.8076
$ws_CUL.Range("I109").Value = 1974.5
$ws_CUL.Range("J109").Value = 5000
$ws_CUL.Range("K109").Value = 5923.5
$ws_CUL.Range("L109").Value = 15000
$ws_CUL.Range("M109").Value = -4883.5
$ws_CUL.Range("N109").Value = -17080
$ws_CUL.Range("H122").Value = 1884
$ws_CUL.Range("I122").Value = 996.3333
$ws_CUL.Range("K122").Value = 8966.9997
$ws_CUL.Range("M122").Value = -6516.9997
$ws_CUL.Range("H123").Value = 4800
$ws_CUL.Range("I123").Value = 2000
$ws_CUL.Range("K123").Value = 6000
$ws_CUL.Range("M123").Value = -3550
$ws_CUL.Range("H129").Value = 1208.9
$ws_CUL.Range("J129").Value = 1561.25
$ws_CUL.Range("L129").Value = 4683.75
$ws_CUL.Range("N129").Value = -14683.75
$ws_CUL.Range("H136").Value = 2348.7058
$ws_CUL.Range("I136").Value = 2009.3334
$ws_CUL.Range("J136").Value = 2421.4285
$ws_CUL.Range("K136").Value = 6028.0002
$ws_CUL.Range("L136").Value = 7264.2855
$ws_CUL.Range("M136").Value = -928.0002000000004
$ws_CUL.Range("N136").Value = -17464.2855
$ws_CUL.Range("H139").Value = 1939
$ws_CUL.Range("I139").Value = 1939
$ws_CUL.Range("K139").Value = 5817
$ws_CUL.Range("M139").Value = -677
$ws_CUL.Range("H140").Value = 4700
$ws_CUL.Range("I140").Value = 4700
$ws_CUL.Range("K140").Value = 14100
$ws_CUL.Range("M140").Value = -8920
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H141").Value = 97999.5
$ws_GSM.Range("J141").Value = 97999.5
$ws_GSM.Range("L141").Value = 97999.5
$ws_GSM.Range("N141").Value = -108359.5
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H81").Value = 5329.4
$ws_WVR.Range("I81").Value = 2882.3333
$ws_WVR.Range("K81").Value = 5764.6666
$ws_WVR.Range("M81").Value = -4703.6666
$ws_WVR.Range("H84").Value = 5329.4
$ws_WVR.Range("I84").Value = 2882.3333
$ws_WVR.Range("K84").Value = 28823.333
$ws_WVR.Range("M84").Value = -23519.333
$ws_WVR.Range("H122").Value = 1186.75
$ws_WVR.Range("I122").Value = 1182.3334
$ws_WVR.Range("K122").Value = 3547.0002
$ws_WVR.Range("M122").Value = -1097.0002
$ws_WVR.Range("H132").Value = 3818.6667
$ws_WVR.Range("I132").Value = 1758.2
$ws_WVR.Range("K132").Value = 5274.6
$ws_WVR.Range("M132").Value = -2744.6
$ws_WVR.Range("H136").Value = 905.8461
$ws_WVR.Range("I136").Value = 907.3333
$ws_WVR.Range("K136").Value = 2721.9999
$ws_WVR.Range("M136").Value = -171.9998999999998
$ws_WVR.Range("H138").Value = 40000
$ws_WVR.Range("J138").Value = 40000
$ws_WVR.Range("L138").Value = 40000
$ws_WVR.Range("N138").Value = -50280
